$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = @("Всё", "хорошо", "прекрасная", "маркиза!", "Все", "хорошо,", "все", "хорошо")
$row3 = @("Daily", "Readings", "to", "make", "your", "Life", "the", "best", "it", "can", "be")

for ($i = 0; $i -lt $row2.Length; $i++) {
    $ws.Cells.Item(2, $i + 1).Value = $row2[$i]
}

for ($i = 0; $i -lt $row3.Length; $i++) {
    $ws.Cells.Item(3, $i + 1).Value = $row3[$i]
}
